# "fixes for players with real clubs" - add new club rows 77-92 to the Mannschaft sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 77: SSV Erkrath 1919 e.V.
$ws.Range("A77").Value = 77
$ws.Range("B77").Value = "SSV Erkrath 1919 e.V."
$ws.Range("C77").Value = "SSV Erkrath"
$ws.Range("S77").Value = 1
$ws.Range("T77").Value = 1

# Row 78: KSC Tesla 07
$ws.Range("A78").Value = 78
$ws.Range("B78").Value = "KSC Tesla 07"
$ws.Range("C78").Value = "Tesla"
$ws.Range("S78").Value = 1
$ws.Range("T78").Value = 1

# Row 79: NK Croatia 70
$ws.Range("A79").Value = 79
$ws.Range("B79").Value = "NK Croatia 70"
$ws.Range("C79").Value = "Croatia"
$ws.Range("S79").Value = 1
$ws.Range("T79").Value = 1

# Row 80: SV Eintracht Solingen
$ws.Range("A80").Value = 80
$ws.Range("B80").Value = "SV Eintracht Solingen"
$ws.Range("C80").Value = "SV Eintracht"
$ws.Range("S80").Value = 1
$ws.Range("T80").Value = 1

# Row 81: SC Germania Reusrath 1913 e.V.
$ws.Range("A81").Value = 81
$ws.Range("B81").Value = "SC Germania Reusrath 1913 e.V."
$ws.Range("C81").Value = "Reusrath"
$ws.Range("S81").Value = 1
$ws.Range("T81").Value = 1

# Row 82: SV Eintracht Haarbrück-Jakobsberg
$ws.Range("A82").Value = 82
$ws.Range("B82").Value = "SV Eintracht Haarbrück-Jakobsberg"
$ws.Range("C82").Value = "Haarbrück"
$ws.Range("S82").Value = 1
$ws.Range("T82").Value = 1

# Row 83: DJK Sportfreunde Gerresheim 1923 e.V.
$ws.Range("A83").Value = 83
$ws.Range("B83").Value = "DJK Sportfreunde Gerresheim 1923 e.V."
$ws.Range("C83").Value = "Sportfr. Gerresheim"
$ws.Range("S83").Value = 1
$ws.Range("T83").Value = 1

# Row 84: Spvgg 1904 e.V. Mössingen
$ws.Range("A84").Value = 84
$ws.Range("B84").Value = "Spvgg 1904 e.V. Mössingen"
$ws.Range("C84").Value = "Mössingen"
$ws.Range("S84").Value = 1
$ws.Range("T84").Value = 1

# Row 85: SV 1930 Issum e.V.
$ws.Range("A85").Value = 85
$ws.Range("B85").Value = "SV 1930 Issum e.V."
$ws.Range("C85").Value = "Issum"
$ws.Range("S85").Value = 1
$ws.Range("T85").Value = 1

# Row 86: FSV Gevelsberg e.V.
$ws.Range("A86").Value = 86
$ws.Range("B86").Value = "FSV Gevelsberg e.V."
$ws.Range("C86").Value = "Gevelsberg"
$ws.Range("S86").Value = 1
$ws.Range("T86").Value = 1

# Row 87: Verein für Sport und Freizeit von 1975 Düsseldorf-Süd e.V.
$ws.Range("A87").Value = 87
$ws.Range("B87").Value = "Verein für Sport und Freizeit von 1975 Düsseldorf-Süd e.V."
$ws.Range("C87").Value = "SFD75"
$ws.Range("S87").Value = 1
$ws.Range("T87").Value = 1

# Row 88: Sportring Eller 1892 e.V.
$ws.Range("A88").Value = 88
$ws.Range("B88").Value = "Sportring Eller 1892 e.V."
$ws.Range("C88").Value = "Sportring"
$ws.Range("S88").Value = 1
$ws.Range("T88").Value = 1

# Row 89: TuS Homberg 1912 e.V.
$ws.Range("A89").Value = 89
$ws.Range("B89").Value = "TuS Homberg 1912 e.V."
$ws.Range("C89").Value = "Homberg"
$ws.Range("S89").Value = 1
$ws.Range("T89").Value = 1

# Row 90: DJK Sparta Bilk e.V.
$ws.Range("A90").Value = 90
$ws.Range("B90").Value = "DJK Sparta Bilk e.V."
$ws.Range("C90").Value = "Sparta Bilk"
$ws.Range("S90").Value = 1
$ws.Range("T90").Value = 1

# Row 91: TV Angermund 1909 e.V.
$ws.Range("A91").Value = 91
$ws.Range("B91").Value = "TV Angermund 1909 e.V."
$ws.Range("C91").Value = "Angermund"
$ws.Range("S91").Value = 1
$ws.Range("T91").Value = 1

# Row 92: DJK TuSA 06 Düsseldorf e.V.
$ws.Range("A92").Value = 92
$ws.Range("B92").Value = "DJK TuSA 06 Düsseldorf e.V."
$ws.Range("C92").Value = "Tusa"

# Re-fit the (now much wider) club-name column
$ws.Columns.Item(2).ColumnWidth = 52.7

# Scroll the frozen view down to the newly added rows and move the selection
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 56
$win.ScrollColumn = 2
$ws.Range("A93").Select()

# Window state (minimized) seen in the workbook view
$win.WindowState = -4140